$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New diary entry on row 39.
# Write cells in the same order the new shared strings appear in the
# diff (A, F, C, B) so the shared-string table indices line up.
$ws.Range("A39").Value = "1 joulu"
$ws.Range("F39").Value = "Samoissa puuhissa jatketaan. Törmäykset(tai ainakin leikkaukset) tunnistetaan."
$ws.Range("C39").Value = "Jatkoa 29 marras"
$ws.Range("B39").Value = "9.45-11.15"
$ws.Range("G39").Value = 1.5

# Match formatting used by the other rows: wrap-text style on C/F,
# time-format style on B.
$ws.Range("C34").Copy()
$ws.Range("C39").PasteSpecial(-4122)
$ws.Range("F34").Copy()
$ws.Range("F39").PasteSpecial(-4122)
$ws.Range("B34").Copy()
$ws.Range("B39").PasteSpecial(-4122)

# Row height to match the other two-line entries.
$ws.Rows.Item(39).RowHeight = 43.2

# Leave the cursor where the author last left it.
$ws.Range("D39").Select() | Out-Null
